$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F ("right answer") and give it the header "option E".
# This shifts right answer/min/max/section_type/question_type/solution/positive_mark/negative_mark
# one column to the right (F->G ... M->N).
$ws.Range("F1").EntireColumn.Insert()
$ws.Range("F1").Value = "option E"

# After the insert, "section_type" (was column I) now lives in column J. Remove that
# whole column, shifting question_type/solution/positive_mark/negative_mark back left
# (K->J ... N->M).
$ws.Range("J1").EntireColumn.Delete()

# Append the new trailing header in column N.
$ws.Range("N1").Value = "common_data"

# Restore the selection to match the saved view state.
$ws.Range("J25").Select()
